{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// The diff appends three new red (#FF0000) runs right after the existing\n// red \"Specify your GitHub link here:\" run, so the paragraph ends up\n// reading (all in red):\n//   Specify your GitHub link here: https://github.com/shintjoo/repo759/HW02\n//\n// Find the anchor text, then insert the GitHub link text right after it\n// (as its own trailing range) and force its font color to red so it\n// matches the rest of the line, regardless of whatever formatting\n// happened to be active at the insertion point.\n\nconst searchText = \"Specify your GitHub link here:\";\nconst linkText = \" https://github.com/shintjoo/repo759/HW02\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text: \" + searchText);\n}\n\n// Collapse to a caret right after the matched text, then insert the new\n// text there. `insertText(..., \"End\")` returns a Range over just the\n// inserted text, which we then color red to match the label before it.\nconst anchorRange = results.items[0];\nconst insertionPoint = anchorRange.getRange(\"End\");\nconst insertedRange = insertionPoint.insertText(linkText, \"End\");\ninsertedRange.font.color = \"#FF0000\";\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# The diff appends three new red (#FF0000) runs right after the existing\n# red \"Specify your GitHub link here:\" run, so the paragraph ends up\n# reading (all in red):\n#   Specify your GitHub link here: https://github.com/shintjoo/repo759/HW02\n#\n# Locate the anchor text with Find, collapse the found range to its end\n# (the caret right after the colon) and insert the GitHub link text there,\n# then force the inserted text's font color to red so it matches the rest\n# of the line.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Specify your GitHub link here:\"\n$linkText = \" https://github.com/shintjoo/repo759/HW02\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $searchText\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find target text: $searchText\"\n}\n\n# Collapse the matched range to its end point (wdCollapseEnd = 0), right\n# after the colon, then insert the new text there.\n$rng.Collapse(0)\n$rng.InsertAfter($linkText)\n\n# Re-select exactly the text we just inserted and force it red so it\n# matches the rest of the line regardless of inherited formatting.\n$insertedStart = $rng.Start\n$insertedEnd = $insertedStart + $linkText.Length\n$insertedRange = $d.Range($insertedStart, $insertedEnd)\n$insertedRange.Font.Color = 255\n"}
